$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list data (Coin, Link, Price, Volume(1h) columns)
# Numeric-looking Price strings are prefixed with a literal apostrophe so
# Excel stores them as text (matching the source data's text format)
# instead of auto-converting them to floating point numbers.

$ws.Range("D2").Value = "30.728.36"
$ws.Range("E2").Value = "  +2.40%  "

$ws.Range("D3").Value = "1.895.48"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'248.15"

$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").Value = "'0.2969"
$ws.Range("E8").Value = "  +1.56%  "

$ws.Range("D9").Value = "'0.06821"
$ws.Range("E9").Value = "  +2.65%  "

$ws.Range("D10").Value = "1.896.31"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").Value = "'17.29"
$ws.Range("E11").Value = "  +3.31%  "

$ws.Range("D12").Value = "'92.44"
$ws.Range("E12").Value = "  +6.81%  "

$ws.Range("D13").Value = "'0.07267"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "'5.118"
$ws.Range("E14").Value = "  +4.99%  "

$ws.Range("D15").Value = "'0.6805"
$ws.Range("E15").Value = "  +1.64%  "

$ws.Range("D16").Value = "30.708.79"
$ws.Range("E16").Value = "  +2.25%  "

$ws.Range("D17").Value = "'0.000007991"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").Value = "'13.32"
$ws.Range("E18").Value = "  +4.09%  "

$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").Value = "2.140.11"
$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "'4.861"
$ws.Range("E22").Value = "  +1.78%  "

$ws.Range("D23").Value = "'193.57"
$ws.Range("E23").Value = "  +36.37%  "

$ws.Range("D24").Value = "'6.088"
$ws.Range("E24").Value = "  +6.75%  "

$ws.Range("D25").Value = "'9.445"

$ws.Range("D26").Value = "'155.70"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").Value = "'19.27"
$ws.Range("E27").Value = "  +12.55%  "

$ws.Range("D28").Value = "'1.925"
$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("D29").Value = "'1.406"
$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("D30").Value = "'4.365"
$ws.Range("E30").Value = "  +3.97%  "

$ws.Range("D31").Value = "'0.09017"
$ws.Range("E31").Value = "  +2.86%  "

$ws.Range("D32").Value = "'4.043"
$ws.Range("E32").Value = "  +1.69%  "

$ws.Range("E33").Value = "  +2.59%  "

$ws.Range("E34").Value = "  +4.74%  "

$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("D36").Value = "'2.739"
$ws.Range("E36").Value = "  +2.77%  "

$ws.Range("D37").Value = "'0.01865"
$ws.Range("E37").Value = "  +0.74%  "

$ws.Range("D38").Value = "'2.683"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").Value = "'2.169"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").Value = "'0.9439"
$ws.Range("E40").Value = "  +1.12%  "

$ws.Range("D41").Value = "'0.4444"
$ws.Range("E41").Value = "  +4.46%  "

$ws.Range("D42").Value = "'106.48"
$ws.Range("E42").Value = "  +4.04%  "

$ws.Range("D43").Value = "'5.778"
$ws.Range("E43").Value = "  -0.96%  "

$ws.Range("D44").Value = "'1.001"

$ws.Range("D45").Value = "'7.694"
$ws.Range("E45").Value = "  +2.72%  "

$ws.Range("E46").Value = "  +6.66%  "

$ws.Range("D47").Value = "'0.05866"
$ws.Range("E47").Value = "  +3.80%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.440"
$ws.Range("E48").Value = "  +7.53%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.715"
$ws.Range("E49").Value = "  +5.28%  "

$ws.Range("D50").Value = "'0.3962"
$ws.Range("E50").Value = "  +4.82%  "

$ws.Range("D51").Value = "'33.66"
$ws.Range("E51").Value = "  +3.64%  "
